$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the existing "Italy" sheet: market label + new reference code.
# ---------------------------------------------------------------------------
$italy = $wb.Worksheets.Item("Italy")
$italy.Range("B2").Value = "Italy market"
$italy.Range("B4").Value = "NGC-3443/T1916"
$italy.Range("A1:XFD1048576").Select()

# ---------------------------------------------------------------------------
# 2. Add the "Spain" sheet (18-row layout like Italy) by copying Italy,
#    then updating its two data cells.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Copy($null, $lastSheet)
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"
$spain.Range("B2").Value = "Spain market"
$spain.Range("B4").Value = "NGC-3442/T1592"
$spain.Range("A9").Select()

# ---------------------------------------------------------------------------
# 3. Add "Serbia", "Romania" and "Slovakia" (20-row layout like UK), by
#    copying the UK sheet, then updating its two data cells.
# ---------------------------------------------------------------------------
$uk = $wb.Worksheets.Item("UK")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$uk.Copy($null, $lastSheet)
$serbia = $wb.Worksheets.Item($wb.Worksheets.Count)
$serbia.Name = "Serbia"
$serbia.Range("B2").Value = "Serbia market"
$serbia.Range("B4").Value = "NGC-4305/T3495"
$serbia.Range("A10").Select()

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$uk.Copy($null, $lastSheet)
$romania = $wb.Worksheets.Item($wb.Worksheets.Count)
$romania.Name = "Romania"
$romania.Range("B2").Value = "Romania market"
$romania.Range("B4").Value = "NGC-4307/T3541"
$romania.Range("B2").Select()

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$uk.Copy($null, $lastSheet)
$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"
$slovakia.Range("B2").Value = "Slovakia market"
$slovakia.Range("B4").Value = "NGC-4306/T3555"
$slovakia.Range("A5").Select()

# Slovakia is the last-touched / active sheet.
$slovakia.Activate()
